$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $lastPart = $parts[$parts.Length - 1]
            $firstPart = $parts[0]
            if ($lastPart.Equals("System") -and -not $firstPart.Equals("System")) {
                $rest = $parts[0..($parts.Length - 2)]
                $newParts = @("System") + $rest
                $newVal = $newParts -join ", "
                $cell.Value = $newVal
            }
        }
    }
}
